# Rename the header row suffixes: "_old" -> "_FV2210" and "_new" -> "_FV2304"
# (the "diff" header in column K stays untouched), then turn A1:U82 into an
# Excel Table ("Table1") and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Segmentname_FV2210"
$ws.Range("B1").Value = "Segmentgruppe_FV2210"
$ws.Range("C1").Value = "Segment_FV2210"
$ws.Range("D1").Value = "Datenelement_FV2210"
$ws.Range("E1").Value = "Segment ID_FV2210"
$ws.Range("F1").Value = "Code_FV2210"
$ws.Range("G1").Value = "Qualifier_FV2210"
$ws.Range("H1").Value = "Beschreibung_FV2210"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2210"
$ws.Range("J1").Value = "Bedingung_FV2210"
# K1 ("diff") is unchanged.
$ws.Range("L1").Value = "Segmentname_FV2304"
$ws.Range("M1").Value = "Segmentgruppe_FV2304"
$ws.Range("N1").Value = "Segment_FV2304"
$ws.Range("O1").Value = "Datenelement_FV2304"
$ws.Range("P1").Value = "Segment ID_FV2304"
$ws.Range("Q1").Value = "Code_FV2304"
$ws.Range("R1").Value = "Qualifier_FV2304"
$ws.Range("S1").Value = "Beschreibung_FV2304"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2304"
$ws.Range("U1").Value = "Bedingung_FV2304"

# Turn the header + data range into a proper Excel Table, matching the
# original data extent A1:U82.
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U82"), $null, 1)
$tbl.Name = "Table1"

# Freeze the header row (split/freeze above row 2).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
